# Slide 2 ("Strings in CPRL"), Content Placeholder 2, last paragraph:
#   before: "type " | "Name = " | "string[20];"   (3 runs)
#   after : "type Name " | "= string[20];"        (2 runs)
# The overall visible text does not change, only where the run boundaries
# fall (and, as a side effect, which original run's formatting - i.e. the
# presence/absence of the dirty="0" attribute - each resulting run keeps).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$para = $tr.Paragraphs(12, 1)

# Step 1: remove the first run's text ("type "), which shifts the
# remaining two runs ("Name = " and "string[20];") to the start of the
# paragraph.
$firstRun = $para.Characters(1, 5)
$firstRun.Text = ""

# Step 2: the (formerly second) run now occupies characters 1-7 ("Name = ").
# Replace its whole span with "type Name " - since the replacement starts
# inside that run, the new run inherits its formatting (no dirty attr).
$secondRun = $para.Characters(1, 7)
$secondRun.Text = "type Name "

# Step 3: the (formerly third) run now starts right after - grab it by its
# remaining length and prefix it with "= " - since the replacement starts
# inside that run, the new run inherits its formatting (dirty="0").
$remainingLen = $para.Text.Length - 10
$thirdRun = $para.Characters(11, $remainingLen)
$thirdRun.Text = "= " + $thirdRun.Text
